$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data rows (25-33, 35-43) ---
$ws.Range("A25").Value = "Otu00004"
$ws.Range("B25").Value = "Bacteroidetes"
$ws.Range("C25").Value = "Porphyromonadaceae"
$ws.Range("E25").Value = [double]"12.603463963579401"
$ws.Range("G25").Value = [double]"-0.62048909673570196"
$ws.Range("A26").Value = "Otu00097"
$ws.Range("E26").Value = [double]"12.2390926998055"
$ws.Range("A27").Value = "Otu00073"
$ws.Range("E27").Value = [double]"12.0270125031301"
$ws.Range("A28").Value = "Otu00042"
$ws.Range("B28").Value = "Firmicutes"
$ws.Range("C28").Value = "Lachnospiraceae"
$ws.Range("E28").Value = [double]"10.970855663631699"
$ws.Range("G28").Value = [double]"-0.54622065303771095"
$ws.Range("A29").Value = "Otu00005"
$ws.Range("B29").Value = "Bacteroidetes"
$ws.Range("C29").Value = "Bacteroides"
$ws.Range("E29").Value = [double]"10.7539928883744"
$ws.Range("G29").Value = [double]"-0.55658991312276895"
$ws.Range("A30").Value = "Otu00010"
$ws.Range("E30").Value = [double]"10.590521026570499"
$ws.Range("A31").Value = "Otu00014"
$ws.Range("E31").Value = [double]"10.504217422815501"
$ws.Range("A32").Value = "Otu00092"
$ws.Range("E32").Value = [double]"10.3323152256946"
$ws.Range("A33").Value = "Otu00001"
$ws.Range("E33").Value = [double]"9.8008163847131193"
$ws.Range("A35").Value = "Otu00081"
$ws.Range("B35").Value = "Bacteroidetes"
$ws.Range("C35").Value = "Porphyromonadaceae"
$ws.Range("G35").Value = [double]"-0.61834182123088299"
$ws.Range("H35").Value = [double]"1.6606556005146901E-10"
$ws.Range("A36").Value = "Otu00086"
$ws.Range("B36").Value = "Firmicutes"
$ws.Range("C36").Value = "Ruminococcus"
$ws.Range("G36").Value = [double]"-0.58461145886863497"
$ws.Range("H36").Value = [double]"3.4072736120638401E-9"
$ws.Range("A37").Value = "Otu00026"
$ws.Range("B37").Value = "Firmicutes"
$ws.Range("C37").Value = "Lachnospiraceae"
$ws.Range("G37").Value = [double]"-0.57602798245480202"
$ws.Range("H37").Value = [double]"6.4435586788493901E-9"
$ws.Range("A38").Value = "Otu00018"
$ws.Range("B38").Value = "Bacteroidetes"
$ws.Range("C38").Value = "Porphyromonadaceae"
$ws.Range("G38").Value = [double]"-0.560682964582215"
$ws.Range("H38").Value = [double]"2.09693698898655E-8"
$ws.Range("A39").Value = "Otu00038"
$ws.Range("B39").Value = "Firmicutes"
$ws.Range("C39").Value = "Ruminococcaceae"
$ws.Range("G39").Value = [double]"-0.55818769337937801"
$ws.Range("H39").Value = [double]"2.4530082470993301E-8"
$ws.Range("A40").Value = "Otu00012"
$ws.Range("B40").Value = "Bacteroidetes"
$ws.Range("C40").Value = "Porphyromonadaceae"
$ws.Range("G40").Value = [double]"-0.55172532091891902"
$ws.Range("H40").Value = [double]"3.7719492985284103E-8"
$ws.Range("A41").Value = "Otu00033"
$ws.Range("B41").Value = "Firmicutes"
$ws.Range("C41").Value = "Ruminococcaceae"
$ws.Range("G41").Value = [double]"-0.53215242792609996"
$ws.Range("H41").Value = [double]"1.5511480954107601E-7"
$ws.Range("A42").Value = "Otu00009"
$ws.Range("B42").Value = "Bacteroidetes"
$ws.Range("C42").Value = "Porphyromonadaceae"
$ws.Range("G42").Value = [double]"-0.50230339850246797"
$ws.Range("H42").Value = [double]"1.1392277530317E-6"
$ws.Range("A43").Value = "Otu00050"
$ws.Range("B43").Value = "Actinobacteria"
$ws.Range("C43").Value = "Coriobacteriaceae"
$ws.Range("G43").Value = [double]"-0.50154276003387599"
$ws.Range("H43").Value = [double]"1.1600624557794101E-6"

# --- Apply scientific-notation number format to the new p-value column (H35:H43) ---
$ws.Range("H35:H43").NumberFormat = "0.00E+00"

# --- Update view: scroll position and active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("F45").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
